{"js": "const paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst target = paragraphs.items[paragraphs.items.length - 1];\ntarget.insertText(\"aaaaaa\", Word.InsertLocation.start);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$para = $d.Paragraphs.Last\n$rng = $para.Range\n$rng.Collapse(1)  # wdCollapseStart\n$rng.InsertBefore(\"aaaaaa\")\n"}
